# Fix error in slide 16: "高级筛选" (Advanced Filter) was the wrong label for
# this bullet; it should read "颜色筛选" (Color Filter). The original text is
# replaced in place and two additional runs are appended (matching how the
# authoring tool split the corrected text into three runs: "颜" / "色筛" / "选").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$para = $tr.Paragraphs(10)
$run1 = $para.Runs(1)

# Sanity check: make sure we are editing the expected run before mutating it.
if ($run1.Text -eq "高级筛选") {
    $run1.Text = "颜"
    $run2 = $run1.InsertAfter("色筛")
    $run3 = $run2.InsertAfter("选")
}
